$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.772.07"
$ws.Range("E2").Value = "  +0.91%  "
$ws.Range("D3").Value = "1.649.41"
$ws.Range("E3").Value = "  +1.33%  "
$ws.Range("E4").Value = "  +0.53%  "
$ws.Range("D5").Value = "'216.65"
$ws.Range("E5").Value = "  +1.57%  "
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("E9").Value = "  +0.53%  "
$ws.Range("D10").Value = "'19.23"
$ws.Range("E10").Value = "  +2.33%  "
$ws.Range("D11").Value = "'0.0844"
$ws.Range("D12").Value = "1.879.29"
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.19"
$ws.Range("E13").Value = "  +1.29%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.632.92"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("E15").Value = "  +1.77%  "
$ws.Range("D16").Value = "'65.35"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").Value = "26.782.49"
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").Value = "'217.66"
$ws.Range("E19").Value = "  +1.17%  "
$ws.Range("E20").Value = "  +0.50%  "
$ws.Range("B21").Value = "Toncoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D21").Value = "'2.52"
$ws.Range("E21").Value = "  +17.87%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'4.37"
$ws.Range("E22").Value = "  +1.72%  "
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("E24").Value = "  +1.81%  "
$ws.Range("D25").Value = "'147.18"
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("E28").Value = "  +3.94%  "
$ws.Range("E29").Value = "  +1.29%  "
$ws.Range("D30").Value = "'0.0520"
$ws.Range("E30").Value = "  +1.41%  "
$ws.Range("E31").Value = "  +1.53%  "
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("E33").Value = "  +1.41%  "
$ws.Range("D34").Value = "1.281.34"
$ws.Range("E34").Value = "  +3.02%  "
$ws.Range("E35").Value = "  +2.80%  "
$ws.Range("E36").Value = "  +2.80%  "
$ws.Range("E37").Value = "  +2.02%  "
$ws.Range("D38").Value = "'0.539"
$ws.Range("E38").Value = "  +5.67%  "
$ws.Range("D39").Value = "'0.829"
$ws.Range("E39").Value = "  +4.25%  "
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("E41").Value = "  +2.21%  "
$ws.Range("E42").Value = "  -1.13%  "
$ws.Range("E43").Value = "  +2.02%  "
$ws.Range("D44").Value = "1.790.76"
$ws.Range("E44").Value = "  +1.50%  "
$ws.Range("D45").Value = "'92.08"
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("D46").Value = "'59.79"
$ws.Range("E46").Value = "  +8.92%  "
$ws.Range("E47").Value = "  +1.01%  "
$ws.Range("E48").Value = "  +1.51%  "
$ws.Range("E49").Value = "  +1.08%  "
$ws.Range("D50").Value = "'7.77"
$ws.Range("E50").Value = "  +3.34%  "
$ws.Range("E51").Value = "  +1.66%  "
